# Generate Report for Handback
# Updates the localization-status workbook: for the zh-cn and de-de sheets,
# the "bee47aa8-3d7e-4b02-b979-233ccd75de80" row (row 8) gets a resolved
# "Latest Target File" / "Latest Handback File" / "Latest Handback DateTime"
# and an "Error Detail" describing that the handed-back file is not the latest
# version, plus a new hyperlink on the "Latest Target File" cell pointing at
# the same handback markdown file as column A.

$wb = $excel.ActiveWorkbook

$errorMessage = 'The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/oltest/blob/84bc2b75339bc9e41de8ed0a980aee14c9105052/e2e/bee47aa8-3d7e-4b02-b979-233ccd75de80.md, latest: https://github.com/OpenLocalizationTestOrg/oltest/blob/c0e0c487bfba72e341471bc964987ba8ebd763e5/e2e/bee47aa8-3d7e-4b02-b979-233ccd75de80.md.'
$latestHandbackUrl = 'https://github.com/OpenLocalizationTestOrg/oltest/blob/c0e0c487bfba72e341471bc964987ba8ebd763e5/e2e/bee47aa8-3d7e-4b02-b979-233ccd75de80.md'
$handbackDisplay = 'bee47aa8-3d7e-4b02-b979-233ccd75de80.md'

# cornflower blue (RGB 0x64,0x95,0xED) packed the way Excel's Font.Color expects
$hyperlinkColor = 15570276

function Set-HandbackRow($ws, $latestTargetFile, $handbackDateTime) {
    # Latest Target File (column I) - add hyperlink + value + hyperlink styling
    $targetCell = $ws.Cells.Item(8, 9)
    $targetCell.Value = $handbackDisplay
    $targetCell.Font.Underline = 2
    $targetCell.Font.Color = $hyperlinkColor
    $ws.Hyperlinks.Add($targetCell, $latestHandbackUrl, $null, $null, $handbackDisplay)

    # Latest Handback File (column J)
    $ws.Cells.Item(8, 10).Value = $latestTargetFile

    # Latest Handback DateTime (column K)
    $ws.Cells.Item(8, 11).Value = $handbackDateTime

    # Error Detail (column P)
    $ws.Cells.Item(8, 16).Value = $errorMessage

    # Widen the Error Detail column so the message is readable
    $ws.Columns.Item(16).ColumnWidth = 39.17
}

$wsZhCn = $wb.Worksheets.Item("zh-cn")
Set-HandbackRow $wsZhCn 'bee47aa8-3d7e-4b02-b979-233ccd75de80.f73eac9650821db7d2488196411fd5da109a0fd7.zh-cn.xlf' '2016-08-13 02:52:25'

$wsDeDe = $wb.Worksheets.Item("de-de")
Set-HandbackRow $wsDeDe 'bee47aa8-3d7e-4b02-b979-233ccd75de80.f73eac9650821db7d2488196411fd5da109a0fd7.de-de.xlf' '2016-08-13 02:52:34'
